# Auto-generated edit script
# Commit: Add data for 2024-12-24
# Adds one additional day (2024-12-24) of violent crime incident counts
# to the running 2024 year-to-date totals (column K) across the citywide
# summary, the by-neighborhood summary, and each neighborhood's own sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 7601
$ws.Range("K3").Value = 7865
$ws.Range("K4").Value = 1654
$ws.Range("K5").Value = 564
$ws.Range("K6").Value = 8756
$ws.Range("K7").Value = 26440

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K6").Value = 193
$ws.Range("K7").Value = 794
$ws.Range("K8").Value = 1731
$ws.Range("K11").Value = 467
$ws.Range("K15").Value = 271
$ws.Range("K18").Value = 178
$ws.Range("K19").Value = 762
$ws.Range("K20").Value = 648
$ws.Range("K21").Value = 90
$ws.Range("K25").Value = 122
$ws.Range("K27").Value = 254
$ws.Range("K29").Value = 1453
$ws.Range("K33").Value = 1113
$ws.Range("K34").Value = 152
$ws.Range("K36").Value = 340
$ws.Range("K37").Value = 873
$ws.Range("K42").Value = 980
$ws.Range("K51").Value = 343
$ws.Range("K52").Value = 680
$ws.Range("K53").Value = 331
$ws.Range("K54").Value = 521
$ws.Range("K63").Value = 76
$ws.Range("K64").Value = 157
$ws.Range("K65").Value = 618
$ws.Range("K67").Value = 1029
$ws.Range("K69").Value = 61
$ws.Range("K73").Value = 239
$ws.Range("K75").Value = 86
$ws.Range("K76").Value = 364
$ws.Range("K77").Value = 171
$ws.Range("K79").Value = 647
$ws.Range("K80").Value = 100
$ws.Range("K84").Value = 215
$ws.Range("K90").Value = 255
$ws.Range("K91").Value = 315
$ws.Range("K93").Value = 104
$ws.Range("K94").Value = 355
$ws.Range("K95").Value = 438
$ws.Range("K96").Value = 281
$ws.Range("K101").Value = 26440

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K2").Value = 86
$ws.Range("K3").Value = 59
$ws.Range("K7").Value = 281

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 260
$ws.Range("K3").Value = 251
$ws.Range("K6").Value = 222
$ws.Range("K7").Value = 794

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 159
$ws.Range("K7").Value = 467

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 189
$ws.Range("K5").Value = 24
$ws.Range("K7").Value = 680

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 61

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K2").Value = 87
$ws.Range("K7").Value = 331

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 484
$ws.Range("K3").Value = 522
$ws.Range("K6").Value = 580
$ws.Range("K7").Value = 1731

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K5").Value = 32
$ws.Range("K6").Value = 357
$ws.Range("K7").Value = 1113

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 146
$ws.Range("K3").Value = 150
$ws.Range("K7").Value = 438

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 248
$ws.Range("K3").Value = 289
$ws.Range("K4").Value = 39
$ws.Range("K6").Value = 264
$ws.Range("K7").Value = 873

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 201
$ws.Range("K7").Value = 618

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 281
$ws.Range("K6").Value = 294
$ws.Range("K7").Value = 1029

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 215

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K2").Value = 80
$ws.Range("K4").Value = 37
$ws.Range("K7").Value = 521

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 407
$ws.Range("K3").Value = 515
$ws.Range("K6").Value = 428
$ws.Range("K7").Value = 1453

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 220
$ws.Range("K5").Value = 25
$ws.Range("K6").Value = 258
$ws.Range("K7").Value = 762

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 81
$ws.Range("K3").Value = 70
$ws.Range("K7").Value = 364

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K3").Value = 51
$ws.Range("K7").Value = 193

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 287
$ws.Range("K6").Value = 373
$ws.Range("K7").Value = 980

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 83
$ws.Range("K7").Value = 315

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 90

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 212
$ws.Range("K7").Value = 647

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K3").Value = 45
$ws.Range("K7").Value = 157

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 220
$ws.Range("K6").Value = 186
$ws.Range("K7").Value = 648

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 178

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 127
$ws.Range("K3").Value = 107
$ws.Range("K7").Value = 340

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("K6").Value = 40
$ws.Range("K7").Value = 104

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K2").Value = 59
$ws.Range("K7").Value = 152

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K3").Value = 75
$ws.Range("K6").Value = 165
$ws.Range("K7").Value = 355

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K6").Value = 22
$ws.Range("K7").Value = 122

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K2").Value = 102
$ws.Range("K3").Value = 68
$ws.Range("K7").Value = 271

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K3").Value = 63
$ws.Range("K4").Value = 16
$ws.Range("K7").Value = 239

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 68
$ws.Range("K7").Value = 254

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K2").Value = 33
$ws.Range("K7").Value = 86

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K6").Value = 67
$ws.Range("K7").Value = 255

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K4").Value = 37
$ws.Range("K6").Value = 112
$ws.Range("K7").Value = 343

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 171

$ws = $wb.Worksheets.Item("Rush &amp; Division")
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 100
